# Contacts page TC done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the existing hyperlinks up front - column insertion does not
#    shift hyperlink anchors automatically, so we rebuild them later once
#    all the data lives in its final position.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Insert a new column B ("Description") - this shifts the old
#    UserName..Social columns (B..F) one place to the right (C..G).
# ---------------------------------------------------------------------------
$ws.Columns("B:B").Insert()

# ---------------------------------------------------------------------------
# 3. New column B content - a short description of each test case.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value2 = "Description"
$ws.Range("B2").Value2 = "Valid First and Last Name"
$ws.Range("B3").Value2 = "First Name is blank"
$ws.Range("B4").Value2 = "Last name is blank"

# ---------------------------------------------------------------------------
# 4. Fix up the "First Name" header (E1) - it used to read "First  Name"
#    (double space); it is now simply "First Name".
# ---------------------------------------------------------------------------
$ws.Range("E1").Value2 = "First Name"

# ---------------------------------------------------------------------------
# 5. Row 3 ("First Name is blank" test case) - the first name cell is now
#    empty instead of "Shreyas".
# ---------------------------------------------------------------------------
$ws.Range("E3").Value2 = ""

# ---------------------------------------------------------------------------
# 6. Row 4 ("Last name is blank" test case) - the last name cell is now
#    empty instead of "Shingre".
# ---------------------------------------------------------------------------
$ws.Range("F4").Value2 = ""

# ---------------------------------------------------------------------------
# 7. New columns H (Expected Result), I (Actual Result), J (Result).
# ---------------------------------------------------------------------------
$ws.Range("H1").Value2 = "Expected Result"
$ws.Range("I1").Value2 = "Actual Result"
$ws.Range("J1").Value2 = "Result"

$ws.Range("H2").Value2 = "Vaibhav Hatge"
$ws.Range("H3").Value2 = "The field First Name is required"
$ws.Range("H4").Value2 = "The field Last Name is required"

# ---------------------------------------------------------------------------
# 8. Styling: H2:J4 share the same plain bordered look used throughout the
#    rest of the data rows (same as A2, E2, ...).
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("H2:J4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 9. Header-row styling for the 3 new header cells H1:J1 - bold font on a
#    yellow fill (same as the rest of row 1) with a partial thin border
#    (left/right always, top only on H1).
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> none
$ws.Range("I1").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none

$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("I1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Borders.Item(8).LineStyle = 1       # xlEdgeTop -> thin
$excel.CutCopyMode = $false

Write-Output "structure+values done"
